# Update the TPM-derived expression specificity metrics on the active sheet
# (Gdnf-Ret.xlsx), reflecting new TPM values used by the upstream script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.028010333333333
$ws.Range("N2").Value = 3.084031
$ws.Range("O2").Value = 0.04165745457248912
$ws.Range("P2").Value = 0.04165745457248914
$ws.Range("Q2").Value = 0.717465887809
$ws.Range("R2").Value = 6.457192990281
$ws.Range("S2").Value = 0.04165745457248912
$ws.Range("T2").Value = 0.04165745457248914

# Row 3
$ws.Range("O3").Value = 0.9361395479363341
$ws.Range("P3").Value = 0.9361395479363344
$ws.Range("S3").Value = 0.9361395479363341
$ws.Range("T3").Value = 0.9361395479363344

# Row 4
$ws.Range("O4").Value = 0.02220299749117665
$ws.Range("P4").Value = 0.02220299749117666
$ws.Range("S4").Value = 0.02220299749117665
$ws.Range("T4").Value = 0.02220299749117666
